$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5226.364  # ALC!H40: 4857.6665 -> 5226.364
$ws.Cells.Item(40, 10).Value = 7756.2856  # ALC!J40: 6887 -> 7756.2856
$ws.Cells.Item(40, 12).Value = 7756.2856  # ALC!L40: 6887 -> 7756.2856
$ws.Cells.Item(40, 14).Value = -8106.2856  # ALC!N40: -7237 -> -8106.2856
$ws.Cells.Item(135, 8).Value = 1487.8572  # ALC!H135: 1127.6 -> 1487.8572
$ws.Cells.Item(135, 9).Value = 1235.8334  # ALC!I135: 919.55554 -> 1235.8334
$ws.Cells.Item(135, 11).Value = 11122.5006  # ALC!K135: 8275.99986 -> 11122.5006
$ws.Cells.Item(135, 13).Value = -8587.500599999999  # ALC!M135: -5740.99986 -> -8587.500599999999
$ws.Cells.Item(137, 8).Value = 2416.6667  # ALC!H137: 2500 -> 2416.6667
$ws.Cells.Item(137, 9).Value = 2416.6667  # ALC!I137: 2500 -> 2416.6667
$ws.Cells.Item(137, 11).Value = 7250.000100000001  # ALC!K137: 7500 -> 7250.000100000001
$ws.Cells.Item(137, 13).Value = -4700.000100000001  # ALC!M137: -4950 -> -4700.000100000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 977  # ARM!H17: 951.5 -> 977
$ws.Cells.Item(17, 10).Value = 977  # ARM!J17: 951.5 -> 977
$ws.Cells.Item(17, 12).Value = 977  # ARM!L17: 951.5 -> 977
$ws.Cells.Item(17, 14).Value = -1323  # ARM!N17: -1297.5 -> -1323
$ws.Cells.Item(26, 8).Value = 5415.5  # ARM!H26: 6785.7144 -> 5415.5
$ws.Cells.Item(26, 10).Value = 7995  # ARM!J26: 11501 -> 7995
$ws.Cells.Item(26, 12).Value = 7995  # ARM!L26: 11501 -> 7995
$ws.Cells.Item(26, 14).Value = -8655  # ARM!N26: -12161 -> -8655
$ws.Cells.Item(29, 8).Value = 2950  # ARM!H29: 4900 -> 2950
$ws.Cells.Item(29, 10).Value = 1000  # ARM!J29: 0 -> 1000
$ws.Cells.Item(29, 12).Value = 1000  # ARM!L29: 0 -> 1000
$ws.Cells.Item(29, 14).Value = -1616  # ARM!N29: None -> -1616
$ws.Cells.Item(124, 8).Value = 27499.5  # ARM!H124: 32500 -> 27499.5
$ws.Cells.Item(124, 10).Value = 27499.5  # ARM!J124: 32500 -> 27499.5
$ws.Cells.Item(124, 12).Value = 27499.5  # ARM!L124: 32500 -> 27499.5
$ws.Cells.Item(124, 14).Value = -37319.5  # ARM!N124: -42320 -> -37319.5
$ws.Cells.Item(125, 8).Value = 99999.5  # ARM!H125: 100000 -> 99999.5
$ws.Cells.Item(125, 10).Value = 99999.5  # ARM!J125: 100000 -> 99999.5
$ws.Cells.Item(125, 12).Value = 99999.5  # ARM!L125: 100000 -> 99999.5
$ws.Cells.Item(125, 14).Value = -109839.5  # ARM!N125: -109840 -> -109839.5
$ws.Cells.Item(141, 8).Value = 0  # ARM!H141: 10000 -> 0
$ws.Cells.Item(141, 10).Value = 0  # ARM!J141: 10000 -> 0
$ws.Cells.Item(141, 12).Value = 0  # ARM!L141: 10000 -> 0
$ws.Cells.Item(141, 14).ClearContents()  # ARM!N141: -20360 -> (removed)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 575  # BSM!H8: 360.25 -> 575
$ws.Cells.Item(8, 10).Value = 1000  # BSM!J8: 430.33334 -> 1000
$ws.Cells.Item(8, 12).Value = 1000  # BSM!L8: 430.33334 -> 1000
$ws.Cells.Item(8, 14).Value = -1280  # BSM!N8: -710.33334 -> -1280
$ws.Cells.Item(14, 8).Value = 444.875  # BSM!H14: 449.875 -> 444.875
$ws.Cells.Item(14, 10).Value = 444.875  # BSM!J14: 449.875 -> 444.875
$ws.Cells.Item(14, 12).Value = 444.875  # BSM!L14: 449.875 -> 444.875
$ws.Cells.Item(14, 14).Value = -788.875  # BSM!N14: -793.875 -> -788.875
$ws.Cells.Item(18, 8).Value = 0  # BSM!H18: 500 -> 0
$ws.Cells.Item(18, 10).Value = 0  # BSM!J18: 500 -> 0
$ws.Cells.Item(18, 12).Value = 0  # BSM!L18: 500 -> 0
$ws.Cells.Item(18, 14).ClearContents()  # BSM!N18: -1558 -> (removed)
$ws.Cells.Item(22, 8).Value = 400  # BSM!H22: 370 -> 400
$ws.Cells.Item(22, 9).Value = 385.7143  # BSM!I22: 355.55554 -> 385.7143
$ws.Cells.Item(22, 11).Value = 385.7143  # BSM!K22: 355.55554 -> 385.7143
$ws.Cells.Item(22, 13).Value = -212.7143  # BSM!M22: -182.55554 -> -212.7143
$ws.Cells.Item(23, 8).Value = 0  # BSM!H23: 492.5 -> 0
$ws.Cells.Item(23, 10).Value = 0  # BSM!J23: 492.5 -> 0
$ws.Cells.Item(23, 12).Value = 0  # BSM!L23: 492.5 -> 0
$ws.Cells.Item(23, 14).ClearContents()  # BSM!N23: -1058.5 -> (removed)
$ws.Cells.Item(134, 8).Value = 3658.3333  # BSM!H134: 3587.4 -> 3658.3333
$ws.Cells.Item(134, 10).Value = 4013.5  # BSM!J134: 4014 -> 4013.5
$ws.Cells.Item(134, 12).Value = 12040.5  # BSM!L134: 12042 -> 12040.5
$ws.Cells.Item(134, 14).Value = -17110.5  # BSM!N134: -17112 -> -17110.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 662.5  # CRP!H12: 605 -> 662.5
$ws.Cells.Item(12, 9).Value = 662.5  # CRP!I12: 605 -> 662.5
$ws.Cells.Item(12, 11).Value = 662.5  # CRP!K12: 605 -> 662.5
$ws.Cells.Item(12, 13).Value = -492.5  # CRP!M12: -435 -> -492.5
$ws.Cells.Item(26, 8).Value = 0  # CRP!H26: 20000 -> 0
$ws.Cells.Item(26, 9).Value = 0  # CRP!I26: 20000 -> 0
$ws.Cells.Item(26, 11).Value = 0  # CRP!K26: 20000 -> 0
$ws.Cells.Item(26, 13).ClearContents()  # CRP!M26: -19713 -> (removed)
$ws.Cells.Item(55, 8).Value = 5000  # CRP!H55: 0 -> 5000
$ws.Cells.Item(55, 9).Value = 5000  # CRP!I55: 0 -> 5000
$ws.Cells.Item(55, 11).Value = 5000  # CRP!K55: 0 -> 5000
$ws.Cells.Item(55, 13).Value = -4685  # CRP!M55: None -> -4685
$ws.Cells.Item(107, 8).Value = 899.4545000000001  # CRP!H107: 780.9231 -> 899.4545000000001
$ws.Cells.Item(107, 9).Value = 699.44446  # CRP!I107: 595.7273 -> 699.44446
$ws.Cells.Item(107, 11).Value = 699.44446  # CRP!K107: 595.7273 -> 699.44446
$ws.Cells.Item(107, 13).Value = 1220.55554  # CRP!M107: 1324.2727 -> 1220.55554
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 2974.7144  # CUL!H68: 3303.8333 -> 2974.7144
$ws.Cells.Item(68, 9).Value = 1766.6666  # CUL!I68: 2150 -> 1766.6666
$ws.Cells.Item(68, 11).Value = 5299.9998  # CUL!K68: 6450 -> 5299.9998
$ws.Cells.Item(68, 13).Value = -4488.9998  # CUL!M68: -5639 -> -4488.9998
$ws.Cells.Item(71, 8).Value = 2974.7144  # CUL!H71: 3303.8333 -> 2974.7144
$ws.Cells.Item(71, 9).Value = 1766.6666  # CUL!I71: 2150 -> 1766.6666
$ws.Cells.Item(71, 11).Value = 15899.9994  # CUL!K71: 19350 -> 15899.9994
$ws.Cells.Item(71, 13).Value = -11843.9994  # CUL!M71: -15294 -> -11843.9994
$ws.Cells.Item(129, 8).Value = 3058.0833  # CUL!H129: 3390 -> 3058.0833
$ws.Cells.Item(129, 9).Value = 2099.5  # CUL!I129: 2666.6667 -> 2099.5
$ws.Cells.Item(129, 10).Value = 3537.375  # CUL!J129: 3700 -> 3537.375
$ws.Cells.Item(129, 11).Value = 6298.5  # CUL!K129: 8000.000100000001 -> 6298.5
$ws.Cells.Item(129, 12).Value = 10612.125  # CUL!L129: 11100 -> 10612.125
$ws.Cells.Item(129, 13).Value = -1298.5  # CUL!M129: -3000.000100000001 -> -1298.5
$ws.Cells.Item(129, 14).Value = -20612.125  # CUL!N129: -21100 -> -20612.125
$ws.Cells.Item(132, 8).Value = 3083.1667  # CUL!H132: 3625 -> 3083.1667
$ws.Cells.Item(132, 9).Value = 1999.6666  # CUL!I132: 2000 -> 1999.6666
$ws.Cells.Item(132, 11).Value = 17996.9994  # CUL!K132: 18000 -> 17996.9994
$ws.Cells.Item(132, 13).Value = -15466.9994  # CUL!M132: -15470 -> -15466.9994
$ws.Cells.Item(137, 8).Value = 15000  # CUL!H137: 5332.6665 -> 15000
$ws.Cells.Item(137, 9).Value = 15000  # CUL!I137: 5332.6665 -> 15000
$ws.Cells.Item(137, 11).Value = 45000  # CUL!K137: 15997.9995 -> 45000
$ws.Cells.Item(137, 13).Value = -39900  # CUL!M137: -10897.9995 -> -39900
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 44  # GSM!H6: 2500 -> 44
$ws.Cells.Item(6, 9).Value = 44  # GSM!I6: 0 -> 44
$ws.Cells.Item(6, 10).Value = 0  # GSM!J6: 2500 -> 0
$ws.Cells.Item(6, 11).Value = 44  # GSM!K6: 0 -> 44
$ws.Cells.Item(6, 12).Value = 0  # GSM!L6: 2500 -> 0
$ws.Cells.Item(6, 13).Value = 69  # GSM!M6: None -> 69
$ws.Cells.Item(6, 14).ClearContents()  # GSM!N6: -2726 -> (removed)
$ws.Cells.Item(16, 8).Value = 44  # GSM!H16: 2500 -> 44
$ws.Cells.Item(16, 9).Value = 44  # GSM!I16: 0 -> 44
$ws.Cells.Item(16, 10).Value = 0  # GSM!J16: 2500 -> 0
$ws.Cells.Item(16, 11).Value = 44  # GSM!K16: 0 -> 44
$ws.Cells.Item(16, 12).Value = 0  # GSM!L16: 2500 -> 0
$ws.Cells.Item(16, 13).Value = 206  # GSM!M16: None -> 206
$ws.Cells.Item(16, 14).ClearContents()  # GSM!N16: -3000 -> (removed)
$ws.Cells.Item(19, 8).Value = 314.5  # GSM!H19: 400.66666 -> 314.5
$ws.Cells.Item(19, 10).Value = 314.5  # GSM!J19: 400.66666 -> 314.5
$ws.Cells.Item(19, 12).Value = 314.5  # GSM!L19: 400.66666 -> 314.5
$ws.Cells.Item(19, 14).Value = -890.5  # GSM!N19: -976.66666 -> -890.5
$ws.Cells.Item(80, 8).Value = 52751.25  # GSM!H80: 101503 -> 52751.25
$ws.Cells.Item(80, 9).Value = 2999.5  # GSM!I80: 3000 -> 2999.5
$ws.Cells.Item(80, 10).Value = 102503  # GSM!J80: 200006 -> 102503
$ws.Cells.Item(80, 11).Value = 2999.5  # GSM!K80: 3000 -> 2999.5
$ws.Cells.Item(80, 12).Value = 102503  # GSM!L80: 200006 -> 102503
$ws.Cells.Item(80, 13).Value = -2001.5  # GSM!M80: -2002 -> -2001.5
$ws.Cells.Item(80, 14).Value = -104499  # GSM!N80: -202002 -> -104499
$ws.Cells.Item(83, 8).Value = 52751.25  # GSM!H83: 101503 -> 52751.25
$ws.Cells.Item(83, 9).Value = 2999.5  # GSM!I83: 3000 -> 2999.5
$ws.Cells.Item(83, 10).Value = 102503  # GSM!J83: 200006 -> 102503
$ws.Cells.Item(83, 11).Value = 14997.5  # GSM!K83: 15000 -> 14997.5
$ws.Cells.Item(83, 12).Value = 512515  # GSM!L83: 1000030 -> 512515
$ws.Cells.Item(83, 13).Value = -10005.5  # GSM!M83: -10008 -> -10005.5
$ws.Cells.Item(83, 14).Value = -522499  # GSM!N83: -1010014 -> -522499
$ws.Cells.Item(122, 8).Value = 1500  # GSM!H122: 0 -> 1500
$ws.Cells.Item(122, 9).Value = 1500  # GSM!I122: 0 -> 1500
$ws.Cells.Item(122, 11).Value = 4500  # GSM!K122: 0 -> 4500
$ws.Cells.Item(122, 13).Value = -2050  # GSM!M122: None -> -2050
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 6492.4  # LTW!H40: 7283 -> 6492.4
$ws.Cells.Item(40, 9).Value = 6492.4  # LTW!I40: 7283 -> 6492.4
$ws.Cells.Item(40, 11).Value = 6492.4  # LTW!K40: 7283 -> 6492.4
$ws.Cells.Item(40, 13).Value = -6356.4  # LTW!M40: -7147 -> -6356.4
$ws.Cells.Item(46, 8).Value = 2763.5908  # LTW!H46: 2780.9048 -> 2763.5908
$ws.Cells.Item(46, 10).Value = 3306.0625  # LTW!J46: 3366.4666 -> 3306.0625
$ws.Cells.Item(46, 12).Value = 3306.0625  # LTW!L46: 3366.4666 -> 3306.0625
$ws.Cells.Item(46, 14).Value = -3682.0625  # LTW!N46: -3742.4666 -> -3682.0625
$ws.Cells.Item(82, 8).Value = 2368.625  # LTW!H82: 2491.6667 -> 2368.625
$ws.Cells.Item(82, 9).Value = 2190  # LTW!I82: 2237.5 -> 2190
$ws.Cells.Item(82, 10).Value = 2666.3333  # LTW!J82: 3000 -> 2666.3333
$ws.Cells.Item(82, 11).Value = 2190  # LTW!K82: 2237.5 -> 2190
$ws.Cells.Item(82, 12).Value = 2666.3333  # LTW!L82: 3000 -> 2666.3333
$ws.Cells.Item(82, 13).Value = -1829  # LTW!M82: -1876.5 -> -1829
$ws.Cells.Item(82, 14).Value = -3388.3333  # LTW!N82: -3722 -> -3388.3333
$ws.Cells.Item(85, 8).Value = 2368.625  # LTW!H85: 2491.6667 -> 2368.625
$ws.Cells.Item(85, 9).Value = 2190  # LTW!I85: 2237.5 -> 2190
$ws.Cells.Item(85, 10).Value = 2666.3333  # LTW!J85: 3000 -> 2666.3333
$ws.Cells.Item(85, 11).Value = 2190  # LTW!K85: 2237.5 -> 2190
$ws.Cells.Item(85, 12).Value = 2666.3333  # LTW!L85: 3000 -> 2666.3333
$ws.Cells.Item(85, 13).Value = -942  # LTW!M85: -989.5 -> -942
$ws.Cells.Item(85, 14).Value = -5162.3333  # LTW!N85: -5496 -> -5162.3333
$ws.Cells.Item(127, 8).Value = 99995  # LTW!H127: 0 -> 99995
$ws.Cells.Item(127, 10).Value = 99995  # LTW!J127: 0 -> 99995
$ws.Cells.Item(127, 12).Value = 99995  # LTW!L127: 0 -> 99995
$ws.Cells.Item(127, 14).Value = -109915  # LTW!N127: None -> -109915
$ws.Cells.Item(132, 8).Value = 5500  # LTW!H132: 7250 -> 5500
$ws.Cells.Item(132, 9).Value = 5500  # LTW!I132: 7250 -> 5500
$ws.Cells.Item(132, 11).Value = 16500  # LTW!K132: 21750 -> 16500
$ws.Cells.Item(132, 13).Value = -13970  # LTW!M132: -19220 -> -13970
$ws.Cells.Item(136, 8).Value = 2914.7144  # LTW!H136: 3217.3333 -> 2914.7144
$ws.Cells.Item(136, 9).Value = 2914.7144  # LTW!I136: 3217.3333 -> 2914.7144
$ws.Cells.Item(136, 11).Value = 8744.143199999999  # LTW!K136: 9651.999899999999 -> 8744.143199999999
$ws.Cells.Item(136, 13).Value = -6194.143199999999  # LTW!M136: -7101.999899999999 -> -6194.143199999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 388307.28  # WVR!H4: 494232.1 -> 388307.28
$ws.Cells.Item(4, 9).Value = 57778.855  # WVR!I4: 80930.39999999999 -> 57778.855
$ws.Cells.Item(4, 10).Value = 718835.7  # WVR!J4: 838650.2 -> 718835.7
$ws.Cells.Item(4, 11).Value = 57778.855  # WVR!K4: 80930.39999999999 -> 57778.855
$ws.Cells.Item(4, 12).Value = 718835.7  # WVR!L4: 838650.2 -> 718835.7
$ws.Cells.Item(4, 13).Value = -57665.855  # WVR!M4: -80817.39999999999 -> -57665.855
$ws.Cells.Item(4, 14).Value = -719061.7  # WVR!N4: -838876.2 -> -719061.7
$ws.Cells.Item(40, 8).Value = 30025  # WVR!H40: 13000 -> 30025
$ws.Cells.Item(40, 9).Value = 30025  # WVR!I40: 0 -> 30025
$ws.Cells.Item(40, 10).Value = 0  # WVR!J40: 13000 -> 0
$ws.Cells.Item(40, 11).Value = 30025  # WVR!K40: 0 -> 30025
$ws.Cells.Item(40, 12).Value = 0  # WVR!L40: 13000 -> 0
$ws.Cells.Item(40, 13).Value = -29876  # WVR!M40: None -> -29876
$ws.Cells.Item(40, 14).ClearContents()  # WVR!N40: -13298 -> (removed)
$ws.Cells.Item(44, 8).Value = 9193.166999999999  # WVR!H44: 25000 -> 9193.166999999999
$ws.Cells.Item(44, 10).Value = 9193.166999999999  # WVR!J44: 25000 -> 9193.166999999999
$ws.Cells.Item(44, 12).Value = 9193.166999999999  # WVR!L44: 25000 -> 9193.166999999999
$ws.Cells.Item(44, 14).Value = -10301.167  # WVR!N44: -26108 -> -10301.167
$ws.Cells.Item(62, 8).Value = 5000  # WVR!H62: 6710 -> 5000
$ws.Cells.Item(62, 9).Value = 5000  # WVR!I62: 6710 -> 5000
$ws.Cells.Item(62, 11).Value = 5000  # WVR!K62: 6710 -> 5000
$ws.Cells.Item(62, 13).Value = -4376  # WVR!M62: -6086 -> -4376
$ws.Cells.Item(65, 8).Value = 5000  # WVR!H65: 6710 -> 5000
$ws.Cells.Item(65, 9).Value = 5000  # WVR!I65: 6710 -> 5000
$ws.Cells.Item(65, 11).Value = 25000  # WVR!K65: 33550 -> 25000
$ws.Cells.Item(65, 13).Value = -21880  # WVR!M65: -30430 -> -21880
$ws.Cells.Item(136, 8).Value = 6749.6665  # WVR!H136: 6499.8 -> 6749.6665
$ws.Cells.Item(136, 10).Value = 7999.5  # WVR!J136: 8000 -> 7999.5
$ws.Cells.Item(136, 12).Value = 23998.5  # WVR!L136: 24000 -> 23998.5
$ws.Cells.Item(136, 14).Value = -29098.5  # WVR!N136: -29100 -> -29098.5
